# Auto-generated from the OOXML diff: bulk refresh of cached market-price
# columns (H/I/J/K/L/M/N) produced by the scheduled pricing runner.
# Only literal values change; no formulas are present in this workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 9
$ws.Range("H9").Value = 156.1579
$ws.Range("I9").Value = 161.17647
$ws.Range("K9").Value = 161.17647
$ws.Range("M9").Value = 7.823530000000005
# row 12
$ws.Range("H12").Value = 302.5
$ws.Range("I12").Value = 302.5
$ws.Range("K12").Value = 302.5
$ws.Range("M12").Value = -132.5
# row 21
$ws.Range("H21").Value = 17664.334
$ws.Range("I21").Value = 12197.2
$ws.Range("K21").Value = 12197.2
$ws.Range("M21").Value = -11729.2
# row 23
$ws.Range("H23").Value = 17664.334
$ws.Range("I23").Value = 12197.2
$ws.Range("K23").Value = 12197.2
$ws.Range("M23").Value = -11963.2
# row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# row 32
$ws.Range("H32").Value = 5666.154
$ws.Range("J32").Value = 5666.154
$ws.Range("L32").Value = 5666.154
$ws.Range("N32").Value = -6318.154
# row 38
$ws.Range("H38").Value = 6735
$ws.Range("I38").Value = 688.5714
$ws.Range("J38").Value = 15200
$ws.Range("K38").Value = 2065.7142
$ws.Range("L38").Value = 45600
$ws.Range("M38").Value = -1693.7142
$ws.Range("N38").Value = -46344
# row 58
$ws.Range("H58").Value = 6414.5713
$ws.Range("I58").Value = 225.75
$ws.Range("J58").Value = 14666.333
$ws.Range("K58").Value = 677.25
$ws.Range("L58").Value = 43998.999
$ws.Range("M58").Value = -527.25
$ws.Range("N58").Value = -44298.999
# row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# row 87
$ws.Range("H87").Value = 82742
$ws.Range("J87").Value = 82742
$ws.Range("L87").Value = 82742
$ws.Range("N87").Value = -85238
# row 90
$ws.Range("H90").Value = 82742
$ws.Range("J90").Value = 82742
$ws.Range("L90").Value = 248226
$ws.Range("N90").Value = -260706
# row 138
$ws.Range("H138").Value = 2132.3767
$ws.Range("J138").Value = 2358.862
$ws.Range("L138").Value = 7076.586
$ws.Range("N138").Value = -17356.586

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 25
$ws.Range("H25").Value = 10285.143
$ws.Range("I25").Value = 7999
$ws.Range("K25").Value = 7999
$ws.Range("M25").Value = -7597
# row 61
$ws.Range("H61").Value = 6237.25
$ws.Range("I61").Value = 6006.5386
$ws.Range("K61").Value = 6006.5386
$ws.Range("M61").Value = -5794.5386
# row 122
$ws.Range("H122").Value = 2102.1785
$ws.Range("I122").Value = 2102.1785
$ws.Range("K122").Value = 6306.5355
$ws.Range("M122").Value = -3856.5355
# row 136
$ws.Range("H136").Value = 6237.25
$ws.Range("I136").Value = 6006.5386
$ws.Range("K136").Value = 18019.6158
$ws.Range("M136").Value = -15469.6158

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 80
$ws.Range("H80").Value = 825.0526
$ws.Range("J80").Value = 1161.4286
$ws.Range("L80").Value = 1161.4286
$ws.Range("N80").Value = -3157.4286
# row 83
$ws.Range("H83").Value = 825.0526
$ws.Range("J83").Value = 1161.4286
$ws.Range("L83").Value = 5807.143
$ws.Range("N83").Value = -15791.143

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 1013.2778
$ws.Range("I16").Value = 840.7
$ws.Range("J16").Value = 1229
$ws.Range("K16").Value = 840.7
$ws.Range("L16").Value = 1229
$ws.Range("M16").Value = -553.7
$ws.Range("N16").Value = -1803
# row 58
$ws.Range("H58").Value = 3406.6924
$ws.Range("I58").Value = 3221.889
$ws.Range("K58").Value = 3221.889
$ws.Range("M58").Value = -3018.889
# row 107
$ws.Range("H107").Value = 478.69232
$ws.Range("I107").Value = 455.57144
$ws.Range("J107").Value = 505.66666
$ws.Range("K107").Value = 455.57144
$ws.Range("L107").Value = 505.66666
$ws.Range("M107").Value = 1464.42856
$ws.Range("N107").Value = -4345.66666
# row 113
$ws.Range("H113").Value = 1013.2778
$ws.Range("I113").Value = 840.7
$ws.Range("J113").Value = 1229
$ws.Range("K113").Value = 840.7
$ws.Range("L113").Value = 1229
$ws.Range("M113").Value = 1329.3
$ws.Range("N113").Value = -5569
# row 132
$ws.Range("H132").Value = 8159.9165
$ws.Range("I132").Value = 5768.8887
$ws.Range("J132").Value = 15333
$ws.Range("K132").Value = 17306.6661
$ws.Range("L132").Value = 45999
$ws.Range("M132").Value = -14776.6661
$ws.Range("N132").Value = -51059
# row 136
$ws.Range("H136").Value = 3406.6924
$ws.Range("I136").Value = 3221.889
$ws.Range("K136").Value = 9665.667000000001
$ws.Range("M136").Value = -7115.667000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 6
$ws.Range("H6").Value = 2168.4443
$ws.Range("I6").Value = 702.3333
$ws.Range("J6").Value = 2901.5
$ws.Range("K6").Value = 2106.9999
$ws.Range("L6").Value = 8704.5
$ws.Range("M6").Value = -1993.9999
$ws.Range("N6").Value = -8930.5
# row 11
$ws.Range("H11").Value = 167.47826
$ws.Range("I11").Value = 84.181816
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 252.545448
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = -112.545448
$ws.Range("N11").Value = -6280
# row 17
$ws.Range("H17").Value = 1243.75
$ws.Range("I17").Value = 325
$ws.Range("J17").Value = 4000
$ws.Range("K17").Value = 975
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = -806
$ws.Range("N17").Value = -12338
# row 34
$ws.Range("H34").Value = 1841.1666
$ws.Range("J34").Value = 2674.25
$ws.Range("L34").Value = 8022.75
$ws.Range("N34").Value = -8190.75
# row 39
$ws.Range("H39").Value = 5434.6665
$ws.Range("J39").Value = 5555.909
$ws.Range("L39").Value = 16667.727
$ws.Range("N39").Value = -17255.727
# row 46
$ws.Range("H46").Value = 1667420.1
$ws.Range("J46").Value = 2000744.2
$ws.Range("L46").Value = 6002232.6
$ws.Range("N46").Value = -6002414.6
# row 55
$ws.Range("H55").Value = 12463.929
$ws.Range("I55").Value = 750
$ws.Range("J55").Value = 14416.25
$ws.Range("K55").Value = 2250
$ws.Range("L55").Value = 43248.75
$ws.Range("M55").Value = -2073
$ws.Range("N55").Value = -43602.75
# row 60
$ws.Range("H60").Value = 2551.25
$ws.Range("I60").Value = 1901.6666
$ws.Range("J60").Value = 4500
$ws.Range("K60").Value = 5704.9998
$ws.Range("L60").Value = 13500
$ws.Range("M60").Value = -5453.9998
$ws.Range("N60").Value = -14002
# row 132
$ws.Range("H132").Value = 1835.1904
$ws.Range("I132").Value = 1540
$ws.Range("K132").Value = 13860
$ws.Range("M132").Value = -11330

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 42
$ws.Range("H42").Value = 70000
$ws.Range("J42").Value = 70000
$ws.Range("L42").Value = 70000
$ws.Range("N42").Value = -70970
# row 115
$ws.Range("H115").Value = 70000
$ws.Range("J115").Value = 70000
$ws.Range("L115").Value = 70000
$ws.Range("N115").Value = -72350
# row 128
$ws.Range("H128").Value = 70331.5
$ws.Range("J128").Value = 70331.5
$ws.Range("L128").Value = 70331.5
$ws.Range("N128").Value = -80291.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 5321.0527
$ws.Range("I40").Value = 5321.0527
$ws.Range("K40").Value = 5321.0527
$ws.Range("M40").Value = -5185.0527
# row 93
$ws.Range("H93").Value = 1960.6875
$ws.Range("I93").Value = 1770.3636
$ws.Range("J93").Value = 2379.4
$ws.Range("K93").Value = 1770.3636
$ws.Range("L93").Value = 2379.4
$ws.Range("M93").Value = -522.3635999999999
$ws.Range("N93").Value = -4875.4
# row 100
$ws.Range("H100").Value = 2201.2222
$ws.Range("I100").Value = 2154.348
$ws.Range("J100").Value = 2470.75
$ws.Range("K100").Value = 2154.348
$ws.Range("L100").Value = 2470.75
$ws.Range("M100").Value = -1613.348
$ws.Range("N100").Value = -3552.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 126
$ws.Range("H126").Value = 2614.7896
$ws.Range("I126").Value = 2695.9412
$ws.Range("K126").Value = 8087.823600000001
$ws.Range("M126").Value = -5617.823600000001

